# Update the "Förändrad" date column (C) from 2023-11-13 (45243) to
# 2023-11-14 (45244) for data rows 2 through 18.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
